$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as plain text so values such as
# "1.001" are not auto-converted to numbers by the smart input parsing.
# We restore the original (default/"Normal") style afterwards so no visible
# style change is introduced.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '22.040.61'
$ws.Range('E2').Value = '  -1.86%  '

$ws.Range('D3').Value = '1.555.23'
$ws.Range('E3').Value = '  -1.14%  '

$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  +0.03%  '

$ws.Range('D6').Value = '286.72'
$ws.Range('E6').Value = '  -0.47%  '

$ws.Range('D7').Value = '0.3767'
$ws.Range('E7').Value = '  +1.33%  '

$ws.Range('D8').Value = '0.3239'
$ws.Range('E8').Value = '  -2.33%  '

$ws.Range('D9').Value = '41.36'
$ws.Range('E9').Value = '  -12.97%  '

$ws.Range('E10').Value = '  -2.09%  '

$ws.Range('D11').Value = '0.07310'
$ws.Range('E11').Value = '  -2.88%  '

$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.03%  '

$ws.Range('D13').Value = '19.64'
$ws.Range('E13').Value = '  -5.42%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').Value = '5.703'
$ws.Range('E14').Value = '  -3.98%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '6.849'
$ws.Range('E15').Value = '  -1.21%  '

$ws.Range('D16').Value = '1.560.19'
$ws.Range('E16').Value = '  -0.87%  '

$ws.Range('E17').Value = '  -3.32%  '

$ws.Range('D18').Value = '0.06642'
$ws.Range('E18').Value = '  -1.29%  '

$ws.Range('D19').Value = '85.07'
$ws.Range('E19').Value = '  -3.53%  '

$ws.Range('D20').Value = '6.472'
$ws.Range('E20').Value = '  +1.06%  '

$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  +0.06%  '

$ws.Range('D22').Value = '15.98'
$ws.Range('E22').Value = '  -3.30%  '

$ws.Range('D23').Value = '11.57'
$ws.Range('E23').Value = '  -3.72%  '

$ws.Range('D24').Value = '22.049.21'
$ws.Range('E24').Value = '  -1.80%  '

$ws.Range('D25').Value = '2.243'
$ws.Range('E25').Value = '  -6.01%  '

$ws.Range('D26').Value = '2.519'
$ws.Range('E26').Value = '  -4.00%  '

$ws.Range('D27').Value = '149.34'
$ws.Range('E27').Value = '  -1.28%  '

$ws.Range('D28').Value = '18.93'
$ws.Range('E28').Value = '  -3.61%  '

$ws.Range('D29').Value = '4.841'
$ws.Range('E29').Value = '  -2.06%  '

$ws.Range('D30').Value = '1.735.01'
$ws.Range('E30').Value = '  -0.58%  '

$ws.Range('D31').Value = '120.21'
$ws.Range('E31').Value = '  -4.03%  '

$ws.Range('E32').Value = '  +1.26%  '

$ws.Range('D33').Value = '5.941'
$ws.Range('E33').Value = '  -2.41%  '

$ws.Range('D34').Value = '9.261'
$ws.Range('E34').Value = '  -6.02%  '

$ws.Range('D35').Value = '0.08102'
$ws.Range('E35').Value = '  -3.08%  '

$ws.Range('D36').Value = '1.599'
$ws.Range('E36').Value = '  -19.69%  '

$ws.Range('D37').Value = '5.231'
$ws.Range('E37').Value = '  -2.15%  '

$ws.Range('D38').Value = '0.02283'
$ws.Range('E38').Value = '  -7.00%  '

$ws.Range('D39').Value = '0.06121'
$ws.Range('E39').Value = '  -4.08%  '

$ws.Range('D40').Value = '0.2113'
$ws.Range('E40').Value = '  -5.44%  '

$ws.Range('D41').Value = '1.213'
$ws.Range('E41').Value = '  -7.20%  '

$ws.Range('D42').Value = '10.89'
$ws.Range('E42').Value = '  -4.63%  '

$ws.Range('D43').Value = '1.001'
$ws.Range('E43').Value = '  +0.07%  '

$ws.Range('D44').Value = '0.5936'

$ws.Range('D45').Value = '13.50'
$ws.Range('E45').Value = '  -3.39%  '

$ws.Range('D46').Value = '3.721'
$ws.Range('E46').Value = '  -1.41%  '

$ws.Range('D47').Value = '0.5733'
$ws.Range('E47').Value = '  -5.91%  '

$ws.Range('D48').Value = '1.943'
$ws.Range('E48').Value = '  -5.20%  '

$ws.Range('D49').Value = '119.65'
$ws.Range('E49').Value = '  -4.32%  '

$ws.Range('D50').Value = '1.157'
$ws.Range('E50').Value = '  -4.32%  '

$ws.Range('D51').Value = '0.06938'
$ws.Range('E51').Value = '  -3.74%  '

# Restore original style (default "Normal") on column D now that the
# text values have been entered, so formatting matches the source file.
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Applied cryptos update"